$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.306.19'
$ws.Range('E2').Value = '  +3.99%  '
$ws.Range('D3').Value = '2.041.38'
$ws.Range('E3').Value = '  +2.69%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.93'
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('E6').Value = '  +1.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '65.71'
$ws.Range('E7').Value = '  +9.64%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +9.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.30'
$ws.Range('E10').Value = '  +2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0817'
$ws.Range('E11').Value = '  +10.28%  '
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.917'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.64'
$ws.Range('E14').Value = '  +24.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.75'
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('D16').Value = '2.340.70'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('E17').Value = '  +5.78%  '
$ws.Range('D18').Value = '2.042.39'
$ws.Range('E18').Value = '  +2.81%  '
$ws.Range('D19').Value = '37.228.14'
$ws.Range('E19').Value = '  +3.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.09'
$ws.Range('E20').Value = '  +2.06%  '
$ws.Range('D21').Value = '0.0₃0908'
$ws.Range('E21').Value = '  +6.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.47'
$ws.Range('E22').Value = '  +5.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.29'
$ws.Range('E23').Value = '  +2.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('E25').Value = '  +0.25%  '
$ws.Range('E26').Value = '  +4.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.03'
$ws.Range('E27').Value = '  +4.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.07'
$ws.Range('E28').Value = '  -2.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.08'
$ws.Range('E29').Value = '  +3.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.127'
$ws.Range('E30').Value = '  +29.10%  '
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.19'
$ws.Range('E32').Value = '  +2.83%  '
$ws.Range('E33').Value = '  +4.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0628'
$ws.Range('E34').Value = '  +4.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.66'
$ws.Range('E35').Value = '  +5.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.37'
$ws.Range('E36').Value = '  +11.42%  '
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  +3.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.01'
$ws.Range('E40').Value = '  +30.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.29'
$ws.Range('E41').Value = '  +5.42%  '
$ws.Range('E42').Value = '  +7.45%  '
$ws.Range('E43').Value = '  +5.39%  '
$ws.Range('E44').Value = '  +4.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.33'
$ws.Range('E45').Value = '  +4.31%  '
$ws.Range('E46').Value = '  +2.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '95.15'
$ws.Range('E47').Value = '  +2.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.80'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').Value = '1.388.30'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.76'
$ws.Range('E51').Value = '  +0.47%  '
